$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns stay formatted as text so values
# like "0.5710" or "57.80" keep their trailing zeros instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.266.42'
$ws.Range("E2").Value = '  -2.65%  '
$ws.Range("D3").Value = '1.703.80'
$ws.Range("E3").Value = '  -2.09%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '223.21'
$ws.Range("E5").Value = '  -2.45%  '
$ws.Range("D6").Value = '0.5295'
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.2648'
$ws.Range("E8").Value = '  -4.50%  '
$ws.Range("D9").Value = '0.06569'
$ws.Range("E9").Value = '  -2.64%  '
$ws.Range("D10").Value = '20.79'
$ws.Range("E10").Value = '  -4.26%  '
$ws.Range("D11").Value = '0.07627'
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").Value = '4.559'
$ws.Range("E12").Value = '  -3.15%  '
$ws.Range("D13").Value = '1.710.33'
$ws.Range("E13").Value = '  -3.85%  '
$ws.Range("D14").Value = '1.936.99'
$ws.Range("E14").Value = '  -2.21%  '
$ws.Range("D15").Value = '0.5710'
$ws.Range("E15").Value = '  -4.83%  '
$ws.Range("D16").Value = '0.0₅8157'
$ws.Range("E16").Value = '  -3.25%  '
$ws.Range("D17").Value = '67.48'
$ws.Range("E17").Value = '  -2.65%  '
$ws.Range("D18").Value = '27.231.46'
$ws.Range("E18").Value = '  -2.72%  '
$ws.Range("D19").Value = '215.11'
$ws.Range("E19").Value = '  -4.77%  '
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("E21").Value = '  -3.96%  '
$ws.Range("E22").Value = '  -4.97%  '
$ws.Range("D23").Value = '5.946'
$ws.Range("E23").Value = '  -4.57%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '141.30'
$ws.Range("E25").Value = '  -3.37%  '
$ws.Range("D26").Value = '1.755'
$ws.Range("E26").Value = '  +6.04%  '
$ws.Range("E27").Value = '  -3.11%  '
$ws.Range("D28").Value = '7.231'
$ws.Range("E28").Value = '  -3.23%  '
$ws.Range("E29").Value = '  -4.76%  '
$ws.Range("D30").Value = '0.05373'
$ws.Range("E30").Value = '  -5.62%  '
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("D32").Value = '3.485'
$ws.Range("E32").Value = '  -6.29%  '
$ws.Range("D33").Value = '3.416'
$ws.Range("E33").Value = '  -3.43%  '
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("D35").Value = '2.871'
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("E36").Value = '  -1.40%  '
$ws.Range("D37").Value = '0.9454'
$ws.Range("E37").Value = '  -3.90%  '
$ws.Range("D38").Value = '0.5834'
$ws.Range("E38").Value = '  -2.12%  '
$ws.Range("D39").Value = '0.01626'
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("D40").Value = '5.855'
$ws.Range("E40").Value = '  -2.57%  '
$ws.Range("D41").Value = '1.003'
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = '1.042.05'
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("D43").Value = '0.8374'
$ws.Range("E43").Value = '  -1.32%  '
$ws.Range("D44").Value = '100.66'
$ws.Range("E44").Value = '  -1.53%  '
$ws.Range("D45").Value = '1.843.52'
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("D46").Value = '0.0₈115'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '57.80'
$ws.Range("E47").Value = '  -3.87%  '
$ws.Range("D48").Value = '0.4505'
$ws.Range("E48").Value = '  +1.69%  '
$ws.Range("D49").Value = '1.001'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").Value = '8.036'
$ws.Range("E50").Value = '  -3.57%  '
$ws.Range("D51").Value = '0.05241'
$ws.Range("E51").Value = '  -1.52%  '
